$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 25.74421623169965
$ws.Range("C2").Value = 16.09076496493991
$ws.Range("D2").Value = 6.062778701963372
$ws.Range("E2").Value = 7.510895132600309
$ws.Range("G2").Value = 3.755202878824358
$ws.Range("I2").Value = 37.22399374323896
$ws.Range("L2").Value = 10.52096299665546
$ws.Range("M2").Value = 19.87223249274445
$ws.Range("N2").Value = 22.48800353801506
$ws.Range("B3").Value = 25.35068316575135
$ws.Range("C3").Value = 15.49535993910399
$ws.Range("D3").Value = 5.959005001155548
$ws.Range("E3").Value = 7.401725695681718
$ws.Range("G3").Value = 3.760560352911885
$ws.Range("I3").Value = 37.09307349668776
$ws.Range("L3").Value = 10.53446577749265
$ws.Range("M3").Value = 19.80497545003726
$ws.Range("N3").Value = 22.5150409596124
$ws.Range("B4").Value = 25.11482910260894
$ws.Range("C4").Value = 15.12403024669351
$ws.Range("D4").Value = 5.896490028350224
$ws.Range("E4").Value = 7.332817514623661
$ws.Range("G4").Value = 3.764014561104385
$ws.Range("I4").Value = 37.02067894786667
$ws.Range("L4").Value = 10.54443621275535
$ws.Range("M4").Value = 19.76868688989041
$ws.Range("N4").Value = 22.53352295997582
$ws.Range("B5").Value = 25.02029332850746
$ws.Range("C5").Value = 14.97153104277883
$ws.Range("D5").Value = 5.871351074616761
$ws.Range("E5").Value = 7.304274556301404
$ws.Range("G5").Value = 3.765463784401814
$ws.Range("I5").Value = 36.99319483773154
$ws.Range("L5").Value = 10.54892115838645
$ws.Range("M5").Value = 19.75516549987156
$ws.Range("N5").Value = 22.54152594055092
$ws.Range("B6").Value = 25.00469452074677
$ws.Range("C6").Value = 14.94614566422019
$ws.Range("D6").Value = 5.867198071159552
$ws.Range("E6").Value = 7.299507391435649
$ws.Range("G6").Value = 3.765706945280221
$ws.Range("I6").Value = 36.9887531426243
$ws.Range("L6").Value = 10.54969134814109
$ws.Range("M6").Value = 19.75299696029321
$ws.Range("N6").Value = 22.54288325619369
$ws.Range("B7").Value = 25.11354761473915
$ws.Range("C7").Value = 15.12197799558009
$ws.Range("D7").Value = 5.896149589001866
$ws.Range("E7").Value = 7.332434432967706
$ws.Range("G7").Value = 3.764033937122956
$ws.Range("I7").Value = 37.02030011106149
$ws.Range("L7").Value = 10.54449499062651
$ws.Range("M7").Value = 19.76849939914139
$ws.Range("N7").Value = 22.53362898426225
$ws.Range("B8").Value = 25.60740240115686
$ws.Range("C8").Value = 15.88681498597672
$ws.Range("D8").Value = 6.026768519995856
$ws.Range("E8").Value = 7.473649021392649
$ws.Range("G8").Value = 3.75701606822645
$ws.Range("I8").Value = 37.17719929919944
$ws.Range("L8").Value = 10.52526988136639
$ws.Range("M8").Value = 19.84800777921816
$ws.Range("N8").Value = 22.49693471891529
$ws.Range("B9").Value = 26.61611010436492
$ws.Range("C9").Value = 17.33053644191845
$ws.Range("D9").Value = 6.290940241805155
$ws.Range("E9").Value = 7.735261523358397
$ws.Range("G9").Value = 3.744552127366935
$ws.Range("I9").Value = 37.547841057315
$ws.Range("L9").Value = 10.50091609420936
$ws.Range("M9").Value = 20.04326924932568
$ws.Range("N9").Value = 22.43997235886866
$ws.Range("B10").Value = 27.37420667513496
$ws.Range("C10").Value = 18.34461137327153
$ws.Range("D10").Value = 6.487898922319217
$ws.Range("E10").Value = 7.91762031994913
$ws.Range("G10").Value = 3.736173983709579
$ws.Range("I10").Value = 37.85780658525731
$ws.Range("L10").Value = 10.49118419378899
$ws.Range("M10").Value = 20.2101038005443
$ws.Range("N10").Value = 22.40736228603944
$ws.Range("B11").Value = 27.7211856491044
$ws.Range("C11").Value = 18.79364858574137
$ws.Range("D11").Value = 6.577704837702053
$ws.Range("E11").Value = 7.998341107015956
$ws.Range("G11").Value = 3.732529094974152
$ws.Range("I11").Value = 38.00679360517035
$ws.Range("L11").Value = 10.48853291993664
$ws.Range("M11").Value = 20.29091804085894
$ws.Range("N11").Value = 22.39455507480939
$ws.Range("B12").Value = 27.85275227628241
$ws.Range("C12").Value = 18.96177078208355
$ws.Range("D12").Value = 6.611708233730114
$ws.Range("E12").Value = 8.028578536834234
$ws.Range("G12").Value = 3.731172589033133
$ws.Range("I12").Value = 38.06433838430628
$ws.Range("L12").Value = 10.48778448377348
$ws.Range("M12").Value = 20.32221229497195
$ws.Range("N12").Value = 22.38999861951777
$ws.Range("B13").Value = 27.82441139480975
$ws.Range("C13").Value = 18.92565022720068
$ws.Range("D13").Value = 6.604385708617768
$ws.Range("E13").Value = 8.022081160698358
$ws.Range("G13").Value = 3.731463684289684
$ws.Range("I13").Value = 38.05189533760593
$ws.Range("L13").Value = 10.487934305388
$ws.Range("M13").Value = 20.31544202158909
$ws.Range("N13").Value = 22.39096685707552
$ws.Range("B14").Value = 27.73200698160717
$ws.Range("C14").Value = 18.80751936109717
$ws.Range("D14").Value = 6.580502591781577
$ws.Range("E14").Value = 8.000835414122324
$ws.Range("G14").Value = 3.732417019830474
$ws.Range("I14").Value = 38.01150539787078
$ws.Range("L14").Value = 10.48846622338166
$ws.Range("M14").Value = 20.29347887146867
$ws.Range("N14").Value = 22.39417432260512
$ws.Range("B15").Value = 27.67542533676466
$ws.Range("C15").Value = 18.73490672063639
$ws.Range("D15").Value = 6.565871960722689
$ws.Range("E15").Value = 7.987778593351457
$ws.Range("G15").Value = 3.733004050652058
$ws.Range("I15").Value = 37.98691145112826
$ws.Range("L15").Value = 10.48882532201491
$ws.Range("M15").Value = 20.2801153896033
$ws.Range("N15").Value = 22.39617724641516
$ws.Range("B16").Value = 27.35156208715942
$ws.Range("C16").Value = 18.315005208286
$ws.Range("D16").Value = 6.482031143229813
$ws.Range("E16").Value = 7.912299323814723
$ws.Range("G16").Value = 3.73641551736671
$ws.Range("I16").Value = 37.84822878783662
$ws.Range("L16").Value = 10.49139321212103
$ws.Range("M16").Value = 20.20492008002874
$ws.Range("N16").Value = 22.40824021818311
$ws.Range("B17").Value = 27.15333026064725
$ws.Range("C17").Value = 18.05415057965463
$ws.Range("D17").Value = 6.430626235484764
$ws.Range("E17").Value = 7.8654166608647
$ws.Range("G17").Value = 3.73855081844337
$ws.Range("I17").Value = 37.76518143971042
$ws.Range("L17").Value = 10.49342351284188
$ws.Range("M17").Value = 20.16003950444234
$ws.Range("N17").Value = 22.41616099797293
$ws.Range("B18").Value = 27.03951948215953
$ws.Range("C18").Value = 17.90296682729386
$ws.Range("D18").Value = 6.401081081080674
$ws.Range("E18").Value = 7.838241082159385
$ws.Range("G18").Value = 3.739794658115768
$ws.Range("I18").Value = 37.71816751822355
$ws.Range("L18").Value = 10.49475843015059
$ws.Range("M18").Value = 20.13468954429123
$ws.Range("N18").Value = 22.42090752711835
$ws.Range("B19").Value = 27.00102470220538
$ws.Range("C19").Value = 17.85158661256949
$ws.Range("D19").Value = 6.391082367209426
$ws.Range("E19").Value = 7.82900410842365
$ws.Range("G19").Value = 3.740218498128564
$ws.Range("I19").Value = 37.70237930496248
$ws.Range("L19").Value = 10.49523911019085
$ws.Range("M19").Value = 20.12618666714145
$ws.Range("N19").Value = 22.42254731839924
$ws.Range("B20").Value = 27.17441192279211
$ws.Range("C20").Value = 18.08203889527938
$ws.Range("D20").Value = 6.436096392448699
$ws.Range("E20").Value = 7.870429187453008
$ws.Range("G20").Value = 3.738321891485984
$ws.Range("I20").Value = 37.77394421564732
$ws.Range("L20").Value = 10.4931900837374
$ws.Range("M20").Value = 20.16476919077426
$ws.Range("N20").Value = 22.41529806682195
$ws.Range("B21").Value = 27.75914472397474
$ws.Range("C21").Value = 18.84227045644869
$ws.Range("D21").Value = 6.587518020577794
$ws.Range("E21").Value = 8.007084819429968
$ws.Range("G21").Value = 3.732136359442817
$ws.Range("I21").Value = 38.02333850448517
$ws.Range("L21").Value = 10.48830304969789
$ws.Range("M21").Value = 20.29991134231054
$ws.Range("N21").Value = 22.39322423586498
$ws.Range("B22").Value = 28.14225297670295
$ws.Range("C22").Value = 19.32788556922589
$ws.Range("D22").Value = 6.686441696534928
$ws.Range("E22").Value = 8.094471412458743
$ws.Range("G22").Value = 3.728232010205743
$ws.Range("I22").Value = 38.19288775271444
$ws.Range("L22").Value = 10.48659862696417
$ws.Range("M22").Value = 20.39225784684667
$ws.Range("N22").Value = 22.38050848588662
$ws.Range("B23").Value = 27.93773618301501
$ws.Range("C23").Value = 19.06977800970495
$ws.Range("D23").Value = 6.633658784729016
$ws.Range("E23").Value = 8.048010394224299
$ws.Range("G23").Value = 3.730303247398366
$ws.Range("I23").Value = 38.10180378645524
$ws.Range("L23").Value = 10.48737197544478
$ws.Range("M23").Value = 20.34260820279565
$ws.Range("N23").Value = 22.38713797430797
$ws.Range("B24").Value = 27.16488040949577
$ws.Range("C24").Value = 18.06943436019825
$ws.Range("D24").Value = 6.433623306304082
$ws.Range("E24").Value = 7.868163715369134
$ws.Range("G24").Value = 3.73842533880667
$ws.Range("I24").Value = 37.76998028583325
$ws.Range("L24").Value = 10.49329509479463
$ws.Range("M24").Value = 20.16262948901749
$ws.Range("N24").Value = 22.4156875977073
$ws.Range("B25").Value = 26.33969700363364
$ws.Range("C25").Value = 16.94735122459267
$ws.Range("D25").Value = 6.218825491111494
$ws.Range("E25").Value = 7.666182077543512
$ws.Range("G25").Value = 3.747786248990258
$ws.Range("I25").Value = 37.44089830721658
$ws.Range("L25").Value = 10.50607280430408
$ws.Range("M25").Value = 19.98629094293643
$ws.Range("N25").Value = 22.45376670094358
